$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 72  # H42: 60.25 -> 72
$ws.Cells.Item(42, 9).Value = 42.666668  # I42: 44 -> 42.666668
$ws.Cells.Item(42, 10).Value = 89.59999999999999  # J42: 109 -> 89.59999999999999
$ws.Cells.Item(42, 11).Value = 128.000004  # K42: 132 -> 128.000004
$ws.Cells.Item(42, 12).Value = 268.8  # L42: 327 -> 268.8
$ws.Cells.Item(42, 13).Value = 101.999996  # M42: 98 -> 101.999996
$ws.Cells.Item(42, 14).Value = -728.8  # N42: -787 -> -728.8
$ws.Cells.Item(64, 8).Value = 83337470  # H64: 90913384 -> 83337470
$ws.Cells.Item(64, 9).Value = 333335070  # I64: 500001340 -> 333335070
$ws.Cells.Item(64, 11).Value = 333335070  # K64: 500001340 -> 333335070
$ws.Cells.Item(64, 13).Value = -333334822  # M64: -500001092 -> -333334822
$ws.Cells.Item(67, 8).Value = 83337470  # H67: 90913384 -> 83337470
$ws.Cells.Item(67, 9).Value = 333335070  # I67: 500001340 -> 333335070
$ws.Cells.Item(67, 11).Value = 333335070  # K67: 500001340 -> 333335070
$ws.Cells.Item(67, 13).Value = -333334212  # M67: -500000482 -> -333334212
$ws.Cells.Item(74, 8).Value = 2795.8333  # H74: 3124.7188 -> 2795.8333
$ws.Cells.Item(74, 9).Value = 2117.1667  # I74: 2167.6667 -> 2117.1667
$ws.Cells.Item(74, 10).Value = 3474.5  # J74: 3499.2173 -> 3474.5
$ws.Cells.Item(74, 11).Value = 2117.1667  # K74: 2167.6667 -> 2117.1667
$ws.Cells.Item(74, 12).Value = 3474.5  # L74: 3499.2173 -> 3474.5
$ws.Cells.Item(74, 13).Value = -1181.1667  # M74: -1231.6667 -> -1181.1667
$ws.Cells.Item(74, 14).Value = -5346.5  # N74: -5371.2173 -> -5346.5
$ws.Cells.Item(77, 8).Value = 2795.8333  # H77: 3124.7188 -> 2795.8333
$ws.Cells.Item(77, 9).Value = 2117.1667  # I77: 2167.6667 -> 2117.1667
$ws.Cells.Item(77, 10).Value = 3474.5  # J77: 3499.2173 -> 3474.5
$ws.Cells.Item(77, 11).Value = 10585.8335  # K77: 10838.3335 -> 10585.8335
$ws.Cells.Item(77, 12).Value = 17372.5  # L77: 17496.0865 -> 17372.5
$ws.Cells.Item(77, 13).Value = -5905.833500000001  # M77: -6158.333500000001 -> -5905.833500000001
$ws.Cells.Item(77, 14).Value = -26732.5  # N77: -26856.0865 -> -26732.5
$ws.Cells.Item(86, 8).Value = 1800.1666  # H86: 2510.3333 -> 1800.1666
$ws.Cells.Item(86, 9).Value = 989.2222  # I86: 1048.25 -> 989.2222
$ws.Cells.Item(86, 10).Value = 2611.111  # J86: 3680 -> 2611.111
$ws.Cells.Item(86, 11).Value = 989.2222  # K86: 1048.25 -> 989.2222
$ws.Cells.Item(86, 12).Value = 2611.111  # L86: 3680 -> 2611.111
$ws.Cells.Item(86, 13).Value = 133.7778  # M86: 74.75 -> 133.7778
$ws.Cells.Item(86, 14).Value = -4857.111  # N86: -5926 -> -4857.111
$ws.Cells.Item(88, 8).Value = 5354.2856  # H88: 4704.0415 -> 5354.2856
$ws.Cells.Item(88, 9).Value = 667.6667  # I88: 600.8570999999999 -> 667.6667
$ws.Cells.Item(88, 10).Value = 7228.933  # J88: 6393.5884 -> 7228.933
$ws.Cells.Item(88, 11).Value = 667.6667  # K88: 600.8570999999999 -> 667.6667
$ws.Cells.Item(88, 12).Value = 7228.933  # L88: 6393.5884 -> 7228.933
$ws.Cells.Item(88, 13).Value = -261.6667  # M88: -194.8570999999999 -> -261.6667
$ws.Cells.Item(88, 14).Value = -8040.933  # N88: -7205.5884 -> -8040.933
$ws.Cells.Item(89, 8).Value = 1800.1666  # H89: 2510.3333 -> 1800.1666
$ws.Cells.Item(89, 9).Value = 989.2222  # I89: 1048.25 -> 989.2222
$ws.Cells.Item(89, 10).Value = 2611.111  # J89: 3680 -> 2611.111
$ws.Cells.Item(89, 11).Value = 4946.111  # K89: 5241.25 -> 4946.111
$ws.Cells.Item(89, 12).Value = 13055.555  # L89: 18400 -> 13055.555
$ws.Cells.Item(89, 13).Value = 669.8890000000001  # M89: 374.75 -> 669.8890000000001
$ws.Cells.Item(89, 14).Value = -24287.555  # N89: -29632 -> -24287.555
$ws.Cells.Item(91, 8).Value = 5354.2856  # H91: 4704.0415 -> 5354.2856
$ws.Cells.Item(91, 9).Value = 667.6667  # I91: 600.8570999999999 -> 667.6667
$ws.Cells.Item(91, 10).Value = 7228.933  # J91: 6393.5884 -> 7228.933
$ws.Cells.Item(91, 11).Value = 667.6667  # K91: 600.8570999999999 -> 667.6667
$ws.Cells.Item(91, 12).Value = 7228.933  # L91: 6393.5884 -> 7228.933
$ws.Cells.Item(91, 13).Value = 736.3333  # M91: 803.1429000000001 -> 736.3333
$ws.Cells.Item(91, 14).Value = -10036.933  # N91: -9201.588400000001 -> -10036.933
$ws.Cells.Item(98, 8).Value = 2084.28  # H98: 1552.0312 -> 2084.28
$ws.Cells.Item(98, 9).Value = 1650.1875  # I98: 1202.3572 -> 1650.1875
$ws.Cells.Item(98, 10).Value = 2856  # J98: 3999.75 -> 2856
$ws.Cells.Item(98, 11).Value = 1650.1875  # K98: 1202.3572 -> 1650.1875
$ws.Cells.Item(98, 12).Value = 2856  # L98: 3999.75 -> 2856
$ws.Cells.Item(98, 13).Value = -152.1875  # M98: 295.6428000000001 -> -152.1875
$ws.Cells.Item(98, 14).Value = -5852  # N98: -6995.75 -> -5852
$ws.Cells.Item(122, 8).Value = 2084.28  # H122: 1552.0312 -> 2084.28
$ws.Cells.Item(122, 9).Value = 1650.1875  # I122: 1202.3572 -> 1650.1875
$ws.Cells.Item(122, 10).Value = 2856  # J122: 3999.75 -> 2856
$ws.Cells.Item(122, 11).Value = 4950.5625  # K122: 3607.0716 -> 4950.5625
$ws.Cells.Item(122, 12).Value = 8568  # L122: 11999.25 -> 8568
$ws.Cells.Item(122, 13).Value = -2500.5625  # M122: -1157.0716 -> -2500.5625
$ws.Cells.Item(122, 14).Value = -13468  # N122: -16899.25 -> -13468

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(9, 8).Value = 14000  # H9: 0 -> 14000
$ws.Cells.Item(9, 10).Value = 14000  # J9: 0 -> 14000
$ws.Cells.Item(9, 12).Value = 14000  # L9: 0 -> 14000
$ws.Cells.Item(9, 14).Value = -14340  # N9: None -> -14340
$ws.Cells.Item(20, 8).Value = 14000  # H20: 0 -> 14000
$ws.Cells.Item(20, 10).Value = 14000  # J20: 0 -> 14000
$ws.Cells.Item(20, 12).Value = 14000  # L20: 0 -> 14000
$ws.Cells.Item(20, 14).Value = -14540  # N20: None -> -14540
$ws.Cells.Item(45, 8).Value = 2450.818  # H45: 2838.3076 -> 2450.818
$ws.Cells.Item(45, 9).Value = 1774.4  # I45: 2524.8 -> 1774.4
$ws.Cells.Item(45, 10).Value = 3014.5  # J45: 3034.25 -> 3014.5
$ws.Cells.Item(45, 11).Value = 1774.4  # K45: 2524.8 -> 1774.4
$ws.Cells.Item(45, 12).Value = 3014.5  # L45: 3034.25 -> 3014.5
$ws.Cells.Item(45, 13).Value = -1397.4  # M45: -2147.8 -> -1397.4
$ws.Cells.Item(45, 14).Value = -3768.5  # N45: -3788.25 -> -3768.5
$ws.Cells.Item(95, 8).Value = 30908  # H95: 30802.666 -> 30908
$ws.Cells.Item(95, 10).Value = 30908  # J95: 30802.666 -> 30908
$ws.Cells.Item(95, 12).Value = 30908  # L95: 30802.666 -> 30908
$ws.Cells.Item(95, 14).Value = -36400  # N95: -36294.666 -> -36400
$ws.Cells.Item(122, 8).Value = 2290.5557  # H122: 1745.5 -> 2290.5557
$ws.Cells.Item(122, 9).Value = 5090  # I122: 1485.3334 -> 5090
$ws.Cells.Item(122, 11).Value = 15270  # K122: 4456.0002 -> 15270
$ws.Cells.Item(122, 13).Value = -12820  # M122: -2006.0002 -> -12820

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 390.64  # H80: 482.05884 -> 390.64
$ws.Cells.Item(80, 9).Value = 360  # I80: 340.6 -> 360
$ws.Cells.Item(80, 10).Value = 402.55554  # J80: 541 -> 402.55554
$ws.Cells.Item(80, 11).Value = 360  # K80: 340.6 -> 360
$ws.Cells.Item(80, 12).Value = 402.55554  # L80: 541 -> 402.55554
$ws.Cells.Item(80, 13).Value = 638  # M80: 657.4 -> 638
$ws.Cells.Item(80, 14).Value = -2398.55554  # N80: -2537 -> -2398.55554
$ws.Cells.Item(83, 8).Value = 390.64  # H83: 482.05884 -> 390.64
$ws.Cells.Item(83, 9).Value = 360  # I83: 340.6 -> 360
$ws.Cells.Item(83, 10).Value = 402.55554  # J83: 541 -> 402.55554
$ws.Cells.Item(83, 11).Value = 1800  # K83: 1703 -> 1800
$ws.Cells.Item(83, 12).Value = 2012.7777  # L83: 2705 -> 2012.7777
$ws.Cells.Item(83, 13).Value = 3192  # M83: 3289 -> 3192
$ws.Cells.Item(83, 14).Value = -11996.7777  # N83: -12689 -> -11996.7777

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(31, 8).Value = 5000  # H31: 3397.5 -> 5000
$ws.Cells.Item(31, 10).Value = 5000  # J31: 3397.5 -> 5000
$ws.Cells.Item(31, 12).Value = 15000  # L31: 10192.5 -> 15000
$ws.Cells.Item(31, 14).Value = -15576  # N31: -10768.5 -> -15576
$ws.Cells.Item(63, 8).Value = 1900  # H63: 2330 -> 1900
$ws.Cells.Item(63, 9).Value = 800  # I63: 990 -> 800
$ws.Cells.Item(63, 11).Value = 2400  # K63: 2970 -> 2400
$ws.Cells.Item(63, 13).Value = -1651  # M63: -2221 -> -1651
$ws.Cells.Item(66, 8).Value = 1900  # H66: 2330 -> 1900
$ws.Cells.Item(66, 9).Value = 800  # I66: 990 -> 800
$ws.Cells.Item(66, 11).Value = 7200  # K66: 8910 -> 7200
$ws.Cells.Item(66, 13).Value = -3456  # M66: -5166 -> -3456
$ws.Cells.Item(76, 8).Value = 1464.6  # H76: 2092.2856 -> 1464.6
$ws.Cells.Item(76, 9).Value = 1080.75  # I76: 882 -> 1080.75
$ws.Cells.Item(76, 11).Value = 3242.25  # K76: 2646 -> 3242.25
$ws.Cells.Item(76, 13).Value = -2859.25  # M76: -2263 -> -2859.25
$ws.Cells.Item(79, 8).Value = 1464.6  # H79: 2092.2856 -> 1464.6
$ws.Cells.Item(79, 9).Value = 1080.75  # I79: 882 -> 1080.75
$ws.Cells.Item(79, 11).Value = 3242.25  # K79: 2646 -> 3242.25
$ws.Cells.Item(79, 13).Value = -1916.25  # M79: -1320 -> -1916.25
$ws.Cells.Item(106, 8).Value = 18266.334  # H106: 999999 -> 18266.334
$ws.Cells.Item(106, 9).Value = 999  # I106: 0 -> 999
$ws.Cells.Item(106, 10).Value = 26900  # J106: 999999 -> 26900
$ws.Cells.Item(106, 11).Value = 2997  # K106: 0 -> 2997
$ws.Cells.Item(106, 12).Value = 80700  # L106: 2999997 -> 80700
$ws.Cells.Item(106, 13).Value = -2051  # M106: None -> -2051
$ws.Cells.Item(106, 14).Value = -82592  # N106: -3001889 -> -82592
$ws.Cells.Item(132, 8).Value = 30585.135  # H132: 28369.975 -> 30585.135
$ws.Cells.Item(132, 9).Value = 661.4286  # I132: 739.2308 -> 661.4286
$ws.Cells.Item(132, 10).Value = 48799.566  # J132: 41673.668 -> 48799.566
$ws.Cells.Item(132, 11).Value = 5952.8574  # K132: 6653.077200000001 -> 5952.8574
$ws.Cells.Item(132, 12).Value = 439196.094  # L132: 375063.012 -> 439196.094
$ws.Cells.Item(132, 13).Value = -3422.8574  # M132: -4123.077200000001 -> -3422.8574
$ws.Cells.Item(132, 14).Value = -444256.094  # N132: -380123.012 -> -444256.094
$ws.Cells.Item(134, 8).Value = 47379.52  # H134: 63596.277 -> 47379.52
$ws.Cells.Item(134, 9).Value = 93937.914  # I134: 75248.87 -> 93937.914
$ws.Cells.Item(134, 10).Value = 4402.5386  # J134: 5333.3335 -> 4402.5386
$ws.Cells.Item(134, 11).Value = 281813.742  # K134: 225746.61 -> 281813.742
$ws.Cells.Item(134, 12).Value = 13207.6158  # L134: 16000.0005 -> 13207.6158
$ws.Cells.Item(134, 13).Value = -276743.742  # M134: -220676.61 -> -276743.742
$ws.Cells.Item(134, 14).Value = -23347.6158  # N134: -26140.0005 -> -23347.6158

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 39190  # H15: 32729.666 -> 39190
$ws.Cells.Item(15, 10).Value = 39190  # J15: 32729.666 -> 39190
$ws.Cells.Item(15, 12).Value = 39190  # L15: 32729.666 -> 39190
$ws.Cells.Item(15, 14).Value = -39766  # N15: -33305.666 -> -39766
$ws.Cells.Item(81, 8).Value = 39190  # H81: 32729.666 -> 39190
$ws.Cells.Item(81, 10).Value = 39190  # J81: 32729.666 -> 39190
$ws.Cells.Item(81, 12).Value = 39190  # L81: 32729.666 -> 39190
$ws.Cells.Item(81, 14).Value = -41186  # N81: -34725.666 -> -41186
$ws.Cells.Item(82, 8).Value = 48000  # H82: 48328 -> 48000
$ws.Cells.Item(82, 10).Value = 48000  # J82: 48328 -> 48000
$ws.Cells.Item(82, 12).Value = 48000  # L82: 48328 -> 48000
$ws.Cells.Item(82, 14).Value = -48766  # N82: -49094 -> -48766
$ws.Cells.Item(84, 8).Value = 39190  # H84: 32729.666 -> 39190
$ws.Cells.Item(84, 10).Value = 39190  # J84: 32729.666 -> 39190
$ws.Cells.Item(84, 12).Value = 117570  # L84: 98188.99800000001 -> 117570
$ws.Cells.Item(84, 14).Value = -127554  # N84: -108172.998 -> -127554
$ws.Cells.Item(85, 8).Value = 48000  # H85: 48328 -> 48000
$ws.Cells.Item(85, 10).Value = 48000  # J85: 48328 -> 48000
$ws.Cells.Item(85, 12).Value = 48000  # L85: 48328 -> 48000
$ws.Cells.Item(85, 14).Value = -50652  # N85: -50980 -> -50652
$ws.Cells.Item(102, 8).Value = 4936.1177  # H102: 2456.65 -> 4936.1177
$ws.Cells.Item(102, 9).Value = 5818.8184  # I102: 2624.8572 -> 5818.8184
$ws.Cells.Item(102, 10).Value = 3317.8333  # J102: 2064.1667 -> 3317.8333
$ws.Cells.Item(102, 11).Value = 5818.8184  # K102: 2624.8572 -> 5818.8184
$ws.Cells.Item(102, 12).Value = 3317.8333  # L102: 2064.1667 -> 3317.8333
$ws.Cells.Item(102, 13).Value = -4196.8184  # M102: -1002.8572 -> -4196.8184
$ws.Cells.Item(102, 14).Value = -6561.8333  # N102: -5308.1667 -> -6561.8333
$ws.Cells.Item(122, 8).Value = 2747.7273  # H122: 1406.1482 -> 2747.7273
$ws.Cells.Item(122, 9).Value = 1045.2858  # I122: 647.7273 -> 1045.2858
$ws.Cells.Item(122, 10).Value = 5727  # J122: 4743.2 -> 5727
$ws.Cells.Item(122, 11).Value = 3135.8574  # K122: 1943.1819 -> 3135.8574
$ws.Cells.Item(122, 12).Value = 17181  # L122: 14229.6 -> 17181
$ws.Cells.Item(122, 13).Value = -685.8574000000003  # M122: 506.8181 -> -685.8574000000003
$ws.Cells.Item(122, 14).Value = -22081  # N122: -19129.6 -> -22081
$ws.Cells.Item(126, 8).Value = 1414.2667  # H126: 1213.75 -> 1414.2667
$ws.Cells.Item(126, 9).Value = 1254.5454  # I126: 1084 -> 1254.5454
$ws.Cells.Item(126, 10).Value = 1853.5  # J126: 1603 -> 1853.5
$ws.Cells.Item(126, 11).Value = 3763.6362  # K126: 3252 -> 3763.6362
$ws.Cells.Item(126, 12).Value = 5560.5  # L126: 4809 -> 5560.5
$ws.Cells.Item(126, 13).Value = -1293.6362  # M126: -782 -> -1293.6362
$ws.Cells.Item(126, 14).Value = -10500.5  # N126: -9749 -> -10500.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2251.6  # H7: 2144.125 -> 2251.6
$ws.Cells.Item(7, 9).Value = 2324  # I7: 2063.2727 -> 2324
$ws.Cells.Item(7, 10).Value = 1600  # J7: 2322 -> 1600
$ws.Cells.Item(7, 11).Value = 2324  # K7: 2063.2727 -> 2324
$ws.Cells.Item(7, 12).Value = 1600  # L7: 2322 -> 1600
$ws.Cells.Item(7, 13).Value = -2212  # M7: -1951.2727 -> -2212
$ws.Cells.Item(7, 14).Value = -1824  # N7: -2546 -> -1824
$ws.Cells.Item(40, 8).Value = 2515.3684  # H40: 2267.3333 -> 2515.3684
$ws.Cells.Item(40, 9).Value = 2385.1333  # I40: 2147.762 -> 2385.1333
$ws.Cells.Item(40, 10).Value = 3003.75  # J40: 2685.8333 -> 3003.75
$ws.Cells.Item(40, 11).Value = 2385.1333  # K40: 2147.762 -> 2385.1333
$ws.Cells.Item(40, 12).Value = 3003.75  # L40: 2685.8333 -> 3003.75
$ws.Cells.Item(40, 13).Value = -2249.1333  # M40: -2011.762 -> -2249.1333
$ws.Cells.Item(40, 14).Value = -3275.75  # N40: -2957.8333 -> -3275.75
$ws.Cells.Item(46, 8).Value = 1231.8948  # H46: 1185.3 -> 1231.8948
$ws.Cells.Item(46, 9).Value = 1402.3636  # I46: 1310.5 -> 1402.3636
$ws.Cells.Item(46, 11).Value = 1402.3636  # K46: 1310.5 -> 1402.3636
$ws.Cells.Item(46, 13).Value = -1214.3636  # M46: -1122.5 -> -1214.3636
$ws.Cells.Item(122, 8).Value = 79701.69500000001  # H122: 69374.53 -> 79701.69500000001
$ws.Cells.Item(122, 9).Value = 112956.89  # I122: 101860.8 -> 112956.89
$ws.Cells.Item(122, 10).Value = 4877.5  # J122: 4402 -> 4877.5
$ws.Cells.Item(122, 11).Value = 338870.67  # K122: 305582.4 -> 338870.67
$ws.Cells.Item(122, 12).Value = 14632.5  # L122: 13206 -> 14632.5
$ws.Cells.Item(122, 13).Value = -336420.67  # M122: -303132.4 -> -336420.67
$ws.Cells.Item(122, 14).Value = -19532.5  # N122: -18106 -> -19532.5
$ws.Cells.Item(126, 8).Value = 2251.6  # H126: 2144.125 -> 2251.6
$ws.Cells.Item(126, 9).Value = 2324  # I126: 2063.2727 -> 2324
$ws.Cells.Item(126, 10).Value = 1600  # J126: 2322 -> 1600
$ws.Cells.Item(126, 11).Value = 6972  # K126: 6189.8181 -> 6972
$ws.Cells.Item(126, 12).Value = 4800  # L126: 6966 -> 4800
$ws.Cells.Item(126, 13).Value = -4502  # M126: -3719.8181 -> -4502
$ws.Cells.Item(126, 14).Value = -9740  # N126: -11906 -> -9740

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 0  # H28: 5000 -> 0
$ws.Cells.Item(28, 10).Value = 0  # J28: 5000 -> 0
$ws.Cells.Item(28, 12).Value = 0  # L28: 5000 -> 0
$ws.Cells.Item(28, 14).ClearContents()  # N28: -5696 -> (removed)
$ws.Cells.Item(88, 8).Value = 24090  # H88: 48000 -> 24090
$ws.Cells.Item(88, 9).Value = 5585.5  # I88: 0 -> 5585.5
$ws.Cells.Item(88, 10).Value = 42594.5  # J88: 48000 -> 42594.5
$ws.Cells.Item(88, 11).Value = 5585.5  # K88: 0 -> 5585.5
$ws.Cells.Item(88, 12).Value = 42594.5  # L88: 48000 -> 42594.5
$ws.Cells.Item(88, 13).Value = -5179.5  # M88: None -> -5179.5
$ws.Cells.Item(88, 14).Value = -43406.5  # N88: -48812 -> -43406.5
$ws.Cells.Item(91, 8).Value = 24090  # H91: 48000 -> 24090
$ws.Cells.Item(91, 9).Value = 5585.5  # I91: 0 -> 5585.5
$ws.Cells.Item(91, 10).Value = 42594.5  # J91: 48000 -> 42594.5
$ws.Cells.Item(91, 11).Value = 5585.5  # K91: 0 -> 5585.5
$ws.Cells.Item(91, 12).Value = 42594.5  # L91: 48000 -> 42594.5
$ws.Cells.Item(91, 13).Value = -4181.5  # M91: None -> -4181.5
$ws.Cells.Item(91, 14).Value = -45402.5  # N91: -50808 -> -45402.5
$ws.Cells.Item(122, 8).Value = 18183100  # H122: 33334168 -> 18183100
$ws.Cells.Item(122, 9).Value = 28572414  # I122: 40000640 -> 28572414
$ws.Cells.Item(122, 10).Value = 1797.5  # J122: 1800 -> 1797.5
$ws.Cells.Item(122, 11).Value = 85717242  # K122: 120001920 -> 85717242
$ws.Cells.Item(122, 12).Value = 5392.5  # L122: 5400 -> 5392.5
$ws.Cells.Item(122, 13).Value = -85714792  # M122: -119999470 -> -85714792
$ws.Cells.Item(122, 14).Value = -10292.5  # N122: -10300 -> -10292.5
